# tetris.xlsx - "sets created for bug when fixing holes"
#
# The "visual" sheet (D:M, rows 2-21) renders the tetris board: a "."
# character marks a filled cell, using style 8 (red fill) for the
# currently-falling piece and style 7 (grey fill) for pieces that have
# already landed. Columns O:X / Z hold helper formulas that recompute
# automatically once the board cells change.
#
# This edit reshapes the bottom of the board (rows 19-21) to fix a bug
# in the hole-filling logic: a couple of cells move from "active" (red)
# to "landed" (grey), a new landed piece appears, and two stray marks
# are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("visual")

# --- Row 19: F19 loses its mark; L19 and M19 gain a landed (grey) mark ---
$ws.Range("D19").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").ClearContents() | Out-Null

$ws.Range("AG3").Copy() | Out-Null
$ws.Range("L19:M19").PasteSpecial(-4122) | Out-Null
$ws.Range("L19:M19").Value = "."

# --- Row 20: G20:I20 gain an active (red) mark; J20:M20 gain a landed (grey) mark ---
$ws.Range("F20").Copy() | Out-Null
$ws.Range("G20:I20").PasteSpecial(-4122) | Out-Null
$ws.Range("G20:I20").Value = "."

$ws.Range("AG3").Copy() | Out-Null
$ws.Range("J20:M20").PasteSpecial(-4122) | Out-Null
$ws.Range("J20:M20").Value = "."

# --- Row 21: F21 and G21 lose their mark; H21:M21 gain a landed (grey) mark ---
$ws.Range("D21").Copy() | Out-Null
$ws.Range("F21:G21").PasteSpecial(-4122) | Out-Null
$ws.Range("F21:G21").ClearContents() | Out-Null

$ws.Range("AG3").Copy() | Out-Null
$ws.Range("H21:M21").PasteSpecial(-4122) | Out-Null
$ws.Range("H21:M21").Value = "."

# --- Restore the selection Excel last had before closing the file ---
$ws.Range("I20").Select() | Out-Null
